# Update "想去人数" (F column) figures across the 展览, 演出 and 全部类型 sheets
# to reflect the regenerated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (exhibitions) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 850
$ws.Range("F3").Value = 1423
$ws.Range("F4").Value = 1064
$ws.Range("F5").Value = 494
$ws.Range("F6").Value = 207
$ws.Range("F7").Value = 644
$ws.Range("F8").Value = 219
$ws.Range("F10").Value = 61
$ws.Range("F12").Value = 133
$ws.Range("F13").Value = 1719
$ws.Range("F14").Value = 418
$ws.Range("F15").Value = 37
$ws.Range("F16").Value = 483
$ws.Range("F17").Value = 84
$ws.Range("F19").Value = 111
$ws.Range("F21").Value = 647
$ws.Range("F22").Value = 41
$ws.Range("F23").Value = 229
$ws.Range("F24").Value = 947
$ws.Range("F26").Value = 1502
$ws.Range("F27").Value = 237

# --- 演出 (performances) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 18

# --- 全部类型 (all types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 850
$ws.Range("F4").Value = 1423
$ws.Range("F5").Value = 1064
$ws.Range("F8").Value = 494
$ws.Range("F9").Value = 207
$ws.Range("F10").Value = 644
$ws.Range("F12").Value = 219
$ws.Range("F14").Value = 61
$ws.Range("F16").Value = 133
$ws.Range("F17").Value = 1719
$ws.Range("F19").Value = 418
$ws.Range("F20").Value = 37
$ws.Range("F21").Value = 483
$ws.Range("F22").Value = 84
$ws.Range("F24").Value = 18
$ws.Range("F25").Value = 111
$ws.Range("F29").Value = 647
$ws.Range("F34").Value = 41
$ws.Range("F35").Value = 229
$ws.Range("F36").Value = 947
$ws.Range("F38").Value = 1502
$ws.Range("F39").Value = 237
